# fix NPC HP error
# Insert a new "MAXHP" column right after "SalePrice" (col F) and before
# the old "MAXMP" column (col G), then populate MAXHP with the same value
# as SalePrice for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G; existing G..Y shift right to H..Z.
$salePriceColWidth = $ws.Columns.Item(6).ColumnWidth
$ws.Columns.Item(7).Insert()
# New column inherits SalePrice's width so it renders the same as the source data.
$ws.Columns.Item(7).ColumnWidth = $salePriceColWidth

# New header cell.
$ws.Cells.Item(1, 7).Value2 = "MAXHP"

# Populate MAXHP (col G) with the SalePrice value (col F) for each data row.
$lastRow = 6
for ($r = 2; $r -le $lastRow; $r++) {
    $salePrice = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 7).Value2 = $salePrice
}

# Tidy up the view: the edit also moved the active selection to I8 and
# reset the frozen/scrolled top-left cell back to A1.
$ws.Range("I8").Select()
